# Auto-generated edit script applying the numeric cell changes described in the diff.
# Each sheet's changed cells are set via Range("CellRef").Value = newValue;
# cells that are removed entirely in the diff are cleared via ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 507.1
$ws.Range("I33").Value = 300.57144
$ws.Range("K33").Value = 300.57144
$ws.Range("M33").Value = -71.57144
$ws.Range("J80").Value = 850
$ws.Range("L80").Value = 2550
$ws.Range("N80").Value = -4546
$ws.Range("J83").Value = 850
$ws.Range("L83").Value = 7650
$ws.Range("N83").Value = -17634
$ws.Range("H98").Value = 5188.9316
$ws.Range("I98").Value = 2440.6287
$ws.Range("K98").Value = 2440.6287
$ws.Range("M98").Value = -942.6287000000002
$ws.Range("H112").Value = 3061.625
$ws.Range("J112").Value = 3061.625
$ws.Range("L112").Value = 9184.875
$ws.Range("N112").Value = -11400.875
$ws.Range("H122").Value = 5188.9316
$ws.Range("I122").Value = 2440.6287
$ws.Range("K122").Value = 7321.886100000001
$ws.Range("M122").Value = -4871.886100000001
$ws.Range("H125").Value = 3271.6
$ws.Range("J125").Value = 4186
$ws.Range("L125").Value = 37674
$ws.Range("N125").Value = -42594
$ws.Range("H138").Value = 3452.2163
$ws.Range("J138").Value = 3008.96
$ws.Range("L138").Value = 9026.880000000001
$ws.Range("N138").Value = -19306.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1869.7142
$ws.Range("I2").Value = 2019
$ws.Range("K2").Value = 2019
$ws.Range("M2").Value = -1906
$ws.Range("H32").Value = 2192.2766
$ws.Range("I32").Value = 1966.9546
$ws.Range("K32").Value = 1966.9546
$ws.Range("M32").Value = -1679.9546
$ws.Range("H43").Value = 49792
$ws.Range("I43").Value = 49999
$ws.Range("J43").Value = 49688.5
$ws.Range("K43").Value = 49999
$ws.Range("L43").Value = 49688.5
$ws.Range("M43").Value = -49686
$ws.Range("N43").Value = -50314.5
$ws.Range("H45").Value = 2147.3333
$ws.Range("I45").Value = 1239.3846
$ws.Range("K45").Value = 1239.3846
$ws.Range("M45").Value = -862.3846000000001
$ws.Range("H61").Value = 9687.833000000001
$ws.Range("I61").Value = 6610.7
$ws.Range("J61").Value = 13534.25
$ws.Range("K61").Value = 6610.7
$ws.Range("L61").Value = 13534.25
$ws.Range("M61").Value = -6398.7
$ws.Range("N61").Value = -13958.25
$ws.Range("H74").Value = 3283.1765
$ws.Range("I74").Value = 2370.6155
$ws.Range("K74").Value = 2370.6155
$ws.Range("M74").Value = -1496.6155
$ws.Range("H77").Value = 3283.1765
$ws.Range("I77").Value = 2370.6155
$ws.Range("K77").Value = 11853.0775
$ws.Range("M77").Value = -7485.077499999999
$ws.Range("H102").Value = 1316.675
$ws.Range("I102").Value = 1237.6666
$ws.Range("K102").Value = 1237.6666
$ws.Range("M102").Value = 384.3334
$ws.Range("H116").Value = 1869.7142
$ws.Range("I116").Value = 2019
$ws.Range("K116").Value = 2019
$ws.Range("M116").Value = 275
$ws.Range("H122").Value = 4988.3335
$ws.Range("I122").Value = 3699.2856
$ws.Range("K122").Value = 11097.8568
$ws.Range("M122").Value = -8647.856800000001
$ws.Range("H132").Value = 1407.3334
$ws.Range("I132").Value = 1461
$ws.Range("J132").Value = 817
$ws.Range("K132").Value = 4383
$ws.Range("L132").Value = 2451
$ws.Range("M132").Value = -1853
$ws.Range("N132").Value = -7511
$ws.Range("H136").Value = 9687.833000000001
$ws.Range("I136").Value = 6610.7
$ws.Range("J136").Value = 13534.25
$ws.Range("K136").Value = 19832.1
$ws.Range("L136").Value = 40602.75
$ws.Range("M136").Value = -17282.1
$ws.Range("N136").Value = -45702.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1869.7142
$ws.Range("I3").Value = 2019
$ws.Range("K3").Value = 2019
$ws.Range("M3").Value = -1905
$ws.Range("H86").Value = 12023.467
$ws.Range("I86").Value = 4879.778
$ws.Range("J86").Value = 22739
$ws.Range("K86").Value = 4879.778
$ws.Range("L86").Value = 22739
$ws.Range("M86").Value = -3756.778
$ws.Range("N86").Value = -24985
$ws.Range("H89").Value = 12023.467
$ws.Range("I89").Value = 4879.778
$ws.Range("J89").Value = 22739
$ws.Range("K89").Value = 24398.89
$ws.Range("L89").Value = 113695
$ws.Range("M89").Value = -18782.89
$ws.Range("N89").Value = -124927
$ws.Range("H94").Value = 10322.111
$ws.Range("I94").Value = 4893.467
$ws.Range("K94").Value = 4893.467
$ws.Range("M94").Value = -4442.467
$ws.Range("H99").Value = 3448.1875
$ws.Range("I99").Value = 3440.8572
$ws.Range("K99").Value = 3440.8572
$ws.Range("M99").Value = -1942.8572
$ws.Range("H134").Value = 10187.3125
$ws.Range("I134").Value = 12051.5
$ws.Range("K134").Value = 36154.5
$ws.Range("M134").Value = -33619.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2565
$ws.Range("I31").Value = 1441.8667
$ws.Range("J31").Value = 4249.7
$ws.Range("K31").Value = 1441.8667
$ws.Range("L31").Value = 4249.7
$ws.Range("M31").Value = -1146.8667
$ws.Range("N31").Value = -4839.7
$ws.Range("H34").Value = 2565
$ws.Range("I34").Value = 1441.8667
$ws.Range("J34").Value = 4249.7
$ws.Range("K34").Value = 1441.8667
$ws.Range("L34").Value = 4249.7
$ws.Range("M34").Value = -1239.8667
$ws.Range("N34").Value = -4653.7
$ws.Range("H58").Value = 3535.7407
$ws.Range("I58").Value = 1929
$ws.Range("J58").Value = 5027.7144
$ws.Range("K58").Value = 1929
$ws.Range("L58").Value = 5027.7144
$ws.Range("M58").Value = -1726
$ws.Range("N58").Value = -5433.7144
$ws.Range("H94").Value = 2566.7693
$ws.Range("I94").Value = 2662
$ws.Range("J94").Value = 2507.25
$ws.Range("K94").Value = 2662
$ws.Range("L94").Value = 2507.25
$ws.Range("M94").Value = -2211
$ws.Range("N94").Value = -3409.25
$ws.Range("H99").Value = 3364.5386
$ws.Range("I99").Value = 2914.2856
$ws.Range("K99").Value = 2914.2856
$ws.Range("M99").Value = -1416.2856
$ws.Range("H107").Value = 5661.2
$ws.Range("I107").Value = 6819.909
$ws.Range("J107").Value = 2474.75
$ws.Range("K107").Value = 6819.909
$ws.Range("L107").Value = 2474.75
$ws.Range("M107").Value = -4899.909
$ws.Range("N107").Value = -6314.75
$ws.Range("H126").Value = 3364.5386
$ws.Range("I126").Value = 2914.2856
$ws.Range("K126").Value = 8742.856800000001
$ws.Range("M126").Value = -6272.856800000001
$ws.Range("H132").Value = 3027
$ws.Range("I132").Value = 2464.1428
$ws.Range("K132").Value = 7392.428400000001
$ws.Range("M132").Value = -4862.428400000001
$ws.Range("H136").Value = 3535.7407
$ws.Range("I136").Value = 1929
$ws.Range("J136").Value = 5027.7144
$ws.Range("K136").Value = 5787
$ws.Range("L136").Value = 15083.1432
$ws.Range("M136").Value = -3237
$ws.Range("N136").Value = -20183.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 317.8
$ws.Range("I18").Value = 150
$ws.Range("K18").Value = 450
$ws.Range("M18").Value = -281
$ws.Range("H70").Value = 1999.5
$ws.Range("I70").Value = 1999.5
$ws.Range("K70").Value = 5998.5
$ws.Range("M70").Value = -5683.5
$ws.Range("H73").Value = 1999.5
$ws.Range("I73").Value = 1999.5
$ws.Range("K73").Value = 5998.5
$ws.Range("M73").Value = -4906.5
$ws.Range("H80").Value = 3998.5
$ws.Range("J80").Value = 3998.5
$ws.Range("L80").Value = 11995.5
$ws.Range("N80").Value = -13867.5
$ws.Range("H83").Value = 3998.5
$ws.Range("J83").Value = 3998.5
$ws.Range("L83").Value = 35986.5
$ws.Range("N83").Value = -45346.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 75971.57000000001
$ws.Range("I14").Value = 103760.2
$ws.Range("J14").Value = 6500
$ws.Range("K14").Value = 103760.2
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = -103592.2
$ws.Range("N14").Value = -6836
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2429.5557
$ws.Range("I132").Value = 2455.647
$ws.Range("J132").Value = 1986
$ws.Range("K132").Value = 7366.941
$ws.Range("L132").Value = 5958
$ws.Range("M132").Value = -4836.941
$ws.Range("N132").Value = -11018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2488.2856
$ws.Range("J68").Value = 2666.6667
$ws.Range("L68").Value = 2666.6667
$ws.Range("N68").Value = -4164.6667
$ws.Range("H71").Value = 2488.2856
$ws.Range("J71").Value = 2666.6667
$ws.Range("L71").Value = 13333.3335
$ws.Range("N71").Value = -20821.3335
$ws.Range("H88").Value = 12085
$ws.Range("I88").Value = 10170
$ws.Range("K88").Value = 10170
$ws.Range("M88").Value = -9742
$ws.Range("H91").Value = 12085
$ws.Range("I91").Value = 10170
$ws.Range("K91").Value = 10170
$ws.Range("M91").Value = -8688
$ws.Range("H124").Value = 43666.332
$ws.Range("J124").Value = 43666.332
$ws.Range("L124").Value = 43666.332
$ws.Range("N124").Value = -53486.332
$ws.Range("H132").Value = 4442.381
$ws.Range("I132").Value = 2257.5715
$ws.Range("J132").Value = 5534.7856
$ws.Range("K132").Value = 6772.7145
$ws.Range("L132").Value = 16604.3568
$ws.Range("M132").Value = -4242.7145
$ws.Range("N132").Value = -21664.3568
$ws.Range("H138").Value = 65000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 65000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 65000
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -75280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2036.841
$ws.Range("I132").Value = 1747.2565
$ws.Range("K132").Value = 5241.7695
$ws.Range("M132").Value = -2711.7695
